# Applies the 'Updated cryptos list' price/volume refresh to the crypto table.
# Numeric-looking price strings (column D) are written with a leading single
# quote so Excel keeps them as text (matching the original inline-string cells)
# instead of auto-coercing them into Number values; the cell style is then reset
# back to Normal so no stray number-format/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Looks like a plain number to Excel -> force text storage.
        $ws.Range($cellAddr).Value = "'" + $text
        $ws.Range($cellAddr).Style = "Normal"
    } else {
        $ws.Range($cellAddr).Value = $text
    }
}

# Row 2 (Bitcoin)
Set-TextValue "D2" '38.381.31'
Set-TextValue "E2" '  +1.22%  '

# Row 3 (Ethereum)
Set-TextValue "D3" '2.086.90'
Set-TextValue "E3" '  +1.96%  '

# Row 4 (TetherUSD)
Set-TextValue "E4" '  -0.03%  '

# Row 5 (BNB)
Set-TextValue "D5" '227.83'
Set-TextValue "E5" '  -0.40%  '

# Row 6 (XRP)
Set-TextValue "E6" '  +0.38%  '

# Row 7 (Solana)
Set-TextValue "D7" '60.80'
Set-TextValue "E7" '  -1.14%  '

# Row 8 (USDC)
Set-TextValue "E8" '  +0.00%  '

# Row 9 (Cardano)
Set-TextValue "E9" '  +0.81%  '

# Row 10 (Dogecoin)
Set-TextValue "D10" '0.0836'
Set-TextValue "E10" '  +2.10%  '

# Row 11 (TRON)
Set-TextValue "E11" '  -0.42%  '

# Row 12 (WrappedliquidstakedEther2.0)
Set-TextValue "D12" '2.397.50'
Set-TextValue "E12" '  +2.01%  '

# Row 13 (Chainlink)
Set-TextValue "D13" '14.81'
Set-TextValue "E13" '  +0.87%  '

# Row 14 (Avalanche)
Set-TextValue "D14" '22.25'
Set-TextValue "E14" '  +5.19%  '

# Row 15 (Polygon)
Set-TextValue "D15" '0.783'
Set-TextValue "E15" '  +0.75%  '

# Row 16 (Polkadot)
Set-TextValue "D16" '5.44'
Set-TextValue "E16" '  +4.42%  '

# Row 17 (WrappedEther)
Set-TextValue "D17" '2.079.83'
Set-TextValue "E17" '  +1.23%  '

# Row 18 (WrappedBTC)
Set-TextValue "D18" '38.278.22'
Set-TextValue "E18" '  +1.06%  '

# Row 19 (Uniswap)
Set-TextValue "B19" 'Litecoin'
Set-TextValue "C19" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D19" '70.92'
Set-TextValue "E19" '  +1.82%  '

# Row 20 (Litecoin)
Set-TextValue "B20" 'Uniswap'
Set-TextValue "C20" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D20" '6.05'
Set-TextValue "E20" '  +2.40%  '

# Row 21 (ShibaInu)
Set-TextValue "D21" '0.0₃0831'
Set-TextValue "E21" '  +0.87%  '

# Row 22 (BitcoinCash)
Set-TextValue "D22" '225.29'
Set-TextValue "E22" '  +0.53%  '

# Row 24 (Toncoin)
Set-TextValue "D24" '2.43'
Set-TextValue "E24" '  +0.19%  '

# Row 25 (PancakeSwap)
Set-TextValue "D25" '2.30'
Set-TextValue "E25" '  +1.21%  '

# Row 26 (Monero)
Set-TextValue "D26" '169.45'
Set-TextValue "E26" '  +0.69%  '

# Row 27 (Cosmos)
Set-TextValue "D27" '9.41'
Set-TextValue "E27" '  +0.58%  '

# Row 28 (Kaspa)
Set-TextValue "D28" '0.133'
Set-TextValue "E28" '  +3.09%  '

# Row 29 (EthereumClassic)
Set-TextValue "D29" '19.02'
Set-TextValue "E29" '  +0.86%  '

# Row 30 (ImmutableX)
Set-TextValue "E30" '  +6.47%  '

# Row 31 (Stellar)
Set-TextValue "E31" '  -0.87%  '

# Row 32 (WEMIXToken)
Set-TextValue "E32" '  +5.89%  '

# Row 33 (InternetComputer(DFINITY))
Set-TextValue "D33" '4.78'
Set-TextValue "E33" '  +5.64%  '

# Row 34 (Filecoin)
Set-TextValue "D34" '4.51'
Set-TextValue "E34" '  +2.54%  '

# Row 35 (Hedera)
Set-TextValue "D35" '0.0605'
Set-TextValue "E35" '  +0.14%  '

# Row 36 (LidoDAOToken)
Set-TextValue "D36" '2.38'
Set-TextValue "E36" '  +1.92%  '

# Row 37 (THORChain)
Set-TextValue "D37" '6.40'
Set-TextValue "E37" '  -3.63%  '

# Row 38 (RenderToken)
Set-TextValue "D38" '3.55'
Set-TextValue "E38" '  +3.01%  '

# Row 39 (BinanceUSD)
Set-TextValue "E39" '  +0.02%  '

# Row 40 (InjectiveProtocol)
Set-TextValue "D40" '18.36'
Set-TextValue "E40" '  +0.99%  '

# Row 41 (Maker)
Set-TextValue "D41" '1.538.10'
Set-TextValue "E41" '  -0.42%  '

# Row 42 (Aave)
Set-TextValue "D42" '99.84'
Set-TextValue "E42" '  +3.45%  '

# Row 43 (VeChain)
Set-TextValue "D43" '0.0219'
Set-TextValue "E43" '  +0.68%  '

# Row 44 (Cronos)
Set-TextValue "D44" '0.0937'
Set-TextValue "E44" '  +2.24%  '

# Row 45 (HuobiToken)
Set-TextValue "D45" '2.82'
Set-TextValue "E45" '  +0.52%  '

# Row 46 (FraxShare)
Set-TextValue "D46" '7.79'
Set-TextValue "E46" '  +10.39%  '

# Row 47 (FTXToken)
Set-TextValue "D47" '4.15'
Set-TextValue "E47" '  +0.44%  '

# Row 48 (TrustWalletToken)
Set-TextValue "D48" '1.11'
Set-TextValue "E48" '  +0.23%  '

# Row 49 (ARBITRUM)
Set-TextValue "E49" '  +1.88%  '

# Row 50 (MXToken)
Set-TextValue "E50" '  +0.52%  '

# Row 51 (RocketPoolETH)
Set-TextValue "D51" '2.283.51'
Set-TextValue "E51" '  +1.96%  '
